$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sibirev I. V.")

# Update scores in row 9 (student "Гусев Артём") columns G:J from 2 to 5
$ws.Range("G9").Value = 5
$ws.Range("H9").Value = 5
$ws.Range("I9").Value = 5
$ws.Range("J9").Value = 5

# Mark row 9 as changed, matching O8's "изм" marker but with a trailing space
$ws.Range("O9").Value = "изм "

# Move the active selection to O10, as recorded in the saved view state
$ws.Range("O10").Select()
